$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "first"
$ws.Range("B1").Value = "second"

$ws.Range("A2").Value = "third"
$ws.Range("B2").Value = "four"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 65
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 6
